# Scheduled market-data refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) for the affected leve rows across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 6940.5625
$ws.Range("I9").Value = 12696.125
$ws.Range("J9").Value = 1185
$ws.Range("K9").Value = 12696.125
$ws.Range("L9").Value = 1185
$ws.Range("M9").Value = -12527.125
$ws.Range("N9").Value = -1523
$ws.Range("H18").Value = 1067.125
$ws.Range("I18").Value = 362.42856
$ws.Range("J18").Value = 6000
$ws.Range("K18").Value = 362.42856
$ws.Range("L18").Value = 6000
$ws.Range("M18").Value = -78.42856
$ws.Range("N18").Value = -6568
$ws.Range("H19").Value = 1714.3334
$ws.Range("I19").Value = 1150.3636
$ws.Range("K19").Value = 1150.3636
$ws.Range("M19").Value = -975.3635999999999
$ws.Range("H46").Value = 2266.6667
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 2700
$ws.Range("L46").Value = 15000
$ws.Range("M46").Value = -2581
$ws.Range("N46").Value = -15238
$ws.Range("H60").Value = 2266.6667
$ws.Range("I60").Value = 900
$ws.Range("J60").Value = 5000
$ws.Range("K60").Value = 2700
$ws.Range("L60").Value = 15000
$ws.Range("M60").Value = -2216
$ws.Range("N60").Value = -15968
$ws.Range("H61").Value = 1599
$ws.Range("I61").Value = 1599
$ws.Range("K61").Value = 4797
$ws.Range("M61").Value = -4625
$ws.Range("H69").Value = 3228.5
$ws.Range("I69").Value = 3249.5
$ws.Range("J69").Value = 3207.5
$ws.Range("K69").Value = 9748.5
$ws.Range("L69").Value = 9622.5
$ws.Range("M69").Value = -8874.5
$ws.Range("N69").Value = -11370.5
$ws.Range("H72").Value = 3228.5
$ws.Range("I72").Value = 3249.5
$ws.Range("J72").Value = 3207.5
$ws.Range("K72").Value = 29245.5
$ws.Range("L72").Value = 28867.5
$ws.Range("M72").Value = -24877.5
$ws.Range("N72").Value = -37603.5
$ws.Range("H137").Value = 3227.3125
$ws.Range("I137").Value = 3592.3845
$ws.Range("J137").Value = 1645.3334
$ws.Range("K137").Value = 10777.1535
$ws.Range("L137").Value = 4936.0002
$ws.Range("M137").Value = -8227.1535
$ws.Range("N137").Value = -10036.0002
$ws.Range("H141").Value = 2007.3125
$ws.Range("I141").Value = 1865.5
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 5596.5
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -416.5
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 637.4
$ws.Range("I2").Value = 513.7778
$ws.Range("J2").Value = 1750
$ws.Range("K2").Value = 513.7778
$ws.Range("L2").Value = 1750
$ws.Range("M2").Value = -400.7778
$ws.Range("N2").Value = -1976
$ws.Range("H32").Value = 3611.2812
$ws.Range("I32").Value = 2016.0602
$ws.Range("J32").Value = 13796.154
$ws.Range("K32").Value = 2016.0602
$ws.Range("L32").Value = 13796.154
$ws.Range("M32").Value = -1729.0602
$ws.Range("N32").Value = -14370.154
$ws.Range("H74").Value = 113823.3
$ws.Range("I74").Value = 125970.445
$ws.Range("K74").Value = 125970.445
$ws.Range("M74").Value = -125096.445
$ws.Range("H77").Value = 113823.3
$ws.Range("I77").Value = 125970.445
$ws.Range("K77").Value = 629852.2250000001
$ws.Range("M77").Value = -625484.2250000001
$ws.Range("H110").Value = 26280.867
$ws.Range("I110").Value = 31476.916
$ws.Range("J110").Value = 5496.6665
$ws.Range("K110").Value = 31476.916
$ws.Range("L110").Value = 5496.6665
$ws.Range("M110").Value = -29431.916
$ws.Range("N110").Value = -9586.666499999999
$ws.Range("H116").Value = 637.4
$ws.Range("I116").Value = 513.7778
$ws.Range("J116").Value = 1750
$ws.Range("K116").Value = 513.7778
$ws.Range("L116").Value = 1750
$ws.Range("M116").Value = 1780.2222
$ws.Range("N116").Value = -6338

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 637.4
$ws.Range("I3").Value = 513.7778
$ws.Range("J3").Value = 1750
$ws.Range("K3").Value = 513.7778
$ws.Range("L3").Value = 1750
$ws.Range("M3").Value = -399.7778
$ws.Range("N3").Value = -1978
$ws.Range("H107").Value = 4110
$ws.Range("I107").Value = 4132.2
$ws.Range("J107").Value = 3999
$ws.Range("K107").Value = 4132.2
$ws.Range("L107").Value = 3999
$ws.Range("M107").Value = -2212.2
$ws.Range("N107").Value = -7839

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 247554.5
$ws.Range("J141").Value = 265420.53
$ws.Range("L141").Value = 265420.53
$ws.Range("N141").Value = -275780.53

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 198997.75
$ws.Range("J37").Value = 198997.75
$ws.Range("L37").Value = 596993.25
$ws.Range("N37").Value = -597217.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1677
$ws.Range("I31").Value = 1677
$ws.Range("K31").Value = 1677
$ws.Range("M31").Value = -1385
$ws.Range("H37").Value = 1677
$ws.Range("I37").Value = 1677
$ws.Range("K37").Value = 1677
$ws.Range("M37").Value = -1400
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H80").Value = 4875
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 5500
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 5500
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -7496
$ws.Range("H83").Value = 4875
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 5500
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 27500
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -37484
$ws.Range("H132").Value = 2371.3704
$ws.Range("I132").Value = 1965.5264
$ws.Range("K132").Value = 5896.5792
$ws.Range("M132").Value = -3366.5792

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2660.923
$ws.Range("I68").Value = 2403.818
$ws.Range("K68").Value = 2403.818
$ws.Range("M68").Value = -1654.818
$ws.Range("H71").Value = 2660.923
$ws.Range("I71").Value = 2403.818
$ws.Range("K71").Value = 12019.09
$ws.Range("M71").Value = -8275.09
$ws.Range("H82").Value = 3581.1765
$ws.Range("I82").Value = 2364.8333
$ws.Range("J82").Value = 6500.4
$ws.Range("K82").Value = 2364.8333
$ws.Range("L82").Value = 6500.4
$ws.Range("M82").Value = -2003.8333
$ws.Range("N82").Value = -7222.4
$ws.Range("H85").Value = 3581.1765
$ws.Range("I85").Value = 2364.8333
$ws.Range("J85").Value = 6500.4
$ws.Range("K85").Value = 2364.8333
$ws.Range("L85").Value = 6500.4
$ws.Range("M85").Value = -1116.8333
$ws.Range("N85").Value = -8996.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 32666.666
$ws.Range("I56").Value = 30000
$ws.Range("J56").Value = 34000
$ws.Range("K56").Value = 30000
$ws.Range("L56").Value = 34000
$ws.Range("M56").Value = -29286
$ws.Range("N56").Value = -35428
$ws.Range("H113").Value = 363.83334
$ws.Range("I113").Value = 273.75
$ws.Range("K113").Value = 821.25
$ws.Range("M113").Value = 1348.75
$ws.Range("H122").Value = 73602.36
$ws.Range("I122").Value = 113042.336
$ws.Range("J122").Value = 2610.4
$ws.Range("K122").Value = 339127.008
$ws.Range("L122").Value = 7831.200000000001
$ws.Range("M122").Value = -336677.008
$ws.Range("N122").Value = -12731.2
